# Apply the quarterly database update: drop the oldest quarter column
# (فصل دوم منتهی به 1399/06), shift all quarterly data one column to the
# left, and append data for the newest quarter (فصل چهارم منتهی به 1401/12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the quarter header labels on row 8 and row 24 ---------------
$headers = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

$cols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $headers[$i]
    $ws.Range($cols[$i] + "24").Value = $headers[$i]
}

# --- Update the quarterly data rows --------------------------------------
# Each row is keyed by its row number, values listed for columns E..N
# in the NEW layout (after dropping the oldest quarter and appending the
# newest quarter's figures).
$data = @{
    10 = @(6237, 2931, 4632, 7670, 11146, 16279, 20593, 63529, 70007, 229745)
    11 = @(35697, 64386, 29668, 147876, 151219, -40692, 72684, 172644, 244457, 1862662)
    12 = @(8053, 51148, 14792, 47743, 62155, 163381, 126563, 307280, 240884, 1375296)
    13 = @(40631, 25721, 22405, 11395, 2366, 10079, 12832, 34137, 8841, 89579)
    14 = @(-2555, -14273, 4241, 10728, 808, 26041, 20603, 48160, 44234, -26820)
    15 = @(1184, 627, 405, 1087, 553, 3618, 2060, 3086, 3287, 6190)
    16 = @(6419, 10170, 13597, 21272, 12695, 64712, 28163, 39367, 70909, 121318)
    17 = @(158645, 133821, 155868, 220600, 223289, 296636, 410919, 428763, 542907, 688666)
    19 = @(8763, 8157, 40843, 27830, 85221, 216204, 361426, 601499, 496319, 1025096)
    20 = @(263074, 282688, 286451, 496201, 549452, 756258, 1055843, 1698465, 1721845, 5371732)
    26 = @(186, 216, 291, 368, 265, 279, 515, 505, 484, 572)
    27 = @(412, 367, 328, 423, 402, 577, 1005, 1006, 1110, 1153)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
